# Apply the "Fix(all scripts): Updating routes for __tmp__ files" edit:
# the results table gained 4 new supermarket chains (columns F:I) and one
# new store row (row 16), some store codes were renumbered/split
# (old "Z423" -> "Z423A" plus a brand new "Z423B" row), and the Huff-model
# output figures for the first store were refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): B1:I1 ------------------------------------------
$headers = @(
    "Supermercado",
    "Hipermercado Metro Independencia",
    "Plaza Vea Izaguirre",
    "Tottus Mega Plaza",
    "Plaza Vea Los Olivos",
    "Tottus Los Olivos",
    "Makro Plaza Lima Norte",
    "Makro Comas"
)
$headerCols = @("B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
}

# ---- Store / route labels, column A, rows 2-16 ---------------------------
$storeLabels = @(
    "Z408", "Z651", "Z407", "Z414", "Z409", "Z403", "Z412", "Z405", "Z417",
    "Z423B", "Z411", "Z399", "Z423A", "Z398", "Z396"
)
for ($i = 0; $i -lt $storeLabels.Length; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row).Value = $storeLabels[$i]
}

# ---- Data values, columns B-I, rows 2-16 ----------------------------------
# Row 2 carries the refreshed figures; every other data row is all zeros.
$row2Values = @(8253933.4900000002, 839098.42, 527388.72, 290858.23999999999, 0, 0, 123084.19, 0)
$dataCols = @("B", "C", "D", "E", "F", "G", "H", "I")

for ($r = 2; $r -le 16; $r++) {
    for ($c = 0; $c -lt $dataCols.Length; $c++) {
        if ($r -eq 2) {
            $value = $row2Values[$c]
        } else {
            $value = 0
        }
        $ws.Range($dataCols[$c] + $r).Value = $value
    }
}

# ---- Formatting: reuse the existing header/label style (index 1) ---------
# for the brand-new cells so they match the rest of the table, instead of
# leaving them with the workbook default style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:I1").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---- Update the sheet's selection to span the enlarged table -------------
$ws.Range("A1:I16").Select() | Out-Null
